$wb = $excel.ActiveWorkbook

# Update the "Status" text from "Ready for handoff" to "In Translation"
# wherever it appears across the workbook. The shorter replacement text
# makes the Status column narrower, so the column width that Excel
# auto-fits to also shrinks for every sheet that shows this column.

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2:F3").Value = "In Translation"
$ws1.Columns("E:F").ColumnWidth = 12.5

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2:C3").Value = "In Translation"
$ws2.Columns("C:C").ColumnWidth = 12.5

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2:C3").Value = "In Translation"
$ws3.Columns("C:C").ColumnWidth = 12.5
